# Apply the data refresh for the new export timestamp.
# Updates a handful of "current stock" values that changed between the
# 12:28 PM and 1:04 PM exports, and refreshes the export timestamp shown
# at the bottom of the report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7  - ABIMOL 300MG 5 RECTAL SUPP.      : stock 6:0  -> 7:0
$ws.Range("H7").Value = "7:0"

# Row 9  - CETAL 250MG/5ML 60ML SUSP        : stock 19:0 -> 20:0
$ws.Range("H9").Value = "20:0"

# Row 15 - OPLEX-N SYRUP 125ML              : stock 4:0  -> 5:0
$ws.Range("H15").Value = "5:0"

# Row 17 - TEGRETOL CR 400MG 20 F.C. DIVITABS : stock 0:1 -> 1:0
$ws.Range("H17").Value = "1:0"

# Footer timestamp: 12:28 PM -> 1:04 PM
$ws.Range("A24").Value = "Thursday, 29 May, 2025 1:04 PM"
